$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1594223333333333
$ws.Range("H2").Value = 0.478267
$ws.Range("I2").Value = 0.01552338951653915
$ws.Range("J2").Value = 0.01552338951653915
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 3.160976126143
$ws.Range("R2").Value = 28.448785135287
$ws.Range("S2").Value = 0.002388101281328769
$ws.Range("T2").Value = 0.002388101281328769

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1594223333333333
$ws.Range("H3").Value = 0.478267
$ws.Range("I3").Value = 0.01552338951653915
$ws.Range("J3").Value = 0.01552338951653915
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("Q3").Value = 13.56996699804567
$ws.Range("R3").Value = 122.129702982411
$ws.Range("S3").Value = 0.01025204059834646
$ws.Range("T3").Value = 0.01025204059834646

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1594223333333333
$ws.Range("H4").Value = 0.478267
$ws.Range("I4").Value = 0.01552338951653915
$ws.Range("J4").Value = 0.01552338951653915
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 3.816369522156111
$ws.Range("R4").Value = 34.347325699405
$ws.Range("S4").Value = 0.002883247636863919
$ws.Range("T4").Value = 0.002883247636863919

$ws.Range("I5").Value = 0.1862883666449807
$ws.Range("J5").Value = 0.1862883666449807
$ws.Range("M5").Value = 19.827687
$ws.Range("N5").Value = 59.483061
$ws.Range("O5").Value = 0.1538389073329896
$ws.Range("P5").Value = 0.1538389073329896
$ws.Range("Q5").Value = 37.93327990098901
$ws.Range("R5").Value = 341.399519108901
$ws.Range("S5").Value = 0.02865839877351118
$ws.Range("T5").Value = 0.02865839877351118

$ws.Range("I6").Value = 0.1862883666449807
$ws.Range("J6").Value = 0.1862883666449807
$ws.Range("O6").Value = 0.6604253914664442
$ws.Range("P6").Value = 0.6604253914664441
$ws.Range("S6").Value = 0.1230295674671559
$ws.Range("T6").Value = 0.1230295674671558

$ws.Range("I7").Value = 0.1862883666449807
$ws.Range("J7").Value = 0.1862883666449807
$ws.Range("M7").Value = 23.93873833333333
$ws.Range("N7").Value = 71.816215
$ws.Range("O7").Value = 0.1857357012005663
$ws.Range("P7").Value = 0.1857357012005663
$ws.Range("Q7").Value = 45.79832542620167
$ws.Range("R7").Value = 412.184928835815
$ws.Range("S7").Value = 0.03460040040431368
$ws.Range("T7").Value = 0.03460040040431368

$ws.Range("G8").Value = 8.197245333333333
$ws.Range("H8").Value = 24.591736
$ws.Range("I8").Value = 0.7981882438384801
$ws.Range("J8").Value = 0.7981882438384801
$ws.Range("M8").Value = 19.827687
$ws.Range("N8").Value = 59.483061
$ws.Range("O8").Value = 0.1538389073329896
$ws.Range("P8").Value = 0.1538389073329896
$ws.Range("Q8").Value = 162.532414731544
$ws.Range("R8").Value = 1462.791732583896
$ws.Range("S8").Value = 0.1227924072781497
$ws.Range("T8").Value = 0.1227924072781496

$ws.Range("G9").Value = 8.197245333333333
$ws.Range("H9").Value = 24.591736
$ws.Range("I9").Value = 0.7981882438384801
$ws.Range("J9").Value = 0.7981882438384801
$ws.Range("O9").Value = 0.6604253914664442
$ws.Range("P9").Value = 0.6604253914664441
$ws.Range("Q9").Value = 697.7463340448987
$ws.Range("R9").Value = 6279.717006404087
$ws.Range("S9").Value = 0.5271437834009418
$ws.Range("T9").Value = 0.5271437834009417

$ws.Range("G10").Value = 8.197245333333333
$ws.Range("H10").Value = 24.591736
$ws.Range("I10").Value = 0.7981882438384801
$ws.Range("J10").Value = 0.7981882438384801
$ws.Range("M10").Value = 23.93873833333333
$ws.Range("N10").Value = 71.816215
$ws.Range("O10").Value = 0.1857357012005663
$ws.Range("P10").Value = 0.1857357012005663
$ws.Range("Q10").Value = 196.2317110888044
$ws.Range("R10").Value = 1766.08539979924
$ws.Range("S10").Value = 0.1482520531593887
$ws.Range("T10").Value = 0.1482520531593887
